# MOS-23045: Update Master Data as per 22 April Changes
# Set is_active (column D) to FALSE for the "Others" gender rows
# (OTH / eng, ara, fra) on the master-gender sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("master-gender")

$ws.Range("D4").Value = $false
$ws.Range("D7").Value = $false
$ws.Range("D10").Value = $false

# Update the active cell selection to match the edited workbook state.
$ws.Range("D12").Select()
